$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds the "Förändrad" (last changed) date for each record.
# All data rows (2 through 74) currently store the serial date value 45204
# (2023-10-05) and need to be updated to 45205 (2023-10-06).
$ws.Range("C2:C74").Value = 45205
